$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country (shared string) swaps in column A ---
$ws.Range("A60").Value = "Kuwait"
$ws.Range("A61").Value = "Tailandia"
$ws.Range("A72").Value = "Azerbaiyan"
$ws.Range("A73").Value = "Estonia"
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("A81").Value = "Eslovaquia"
$ws.Range("A114").Value = "Mayotte"
$ws.Range("A115").Value = "Somalia"

# --- Numeric data updates (Casos totales/Nuevos/Activos/Recuperados/Criticos/MuertesHoy/Muertes) ---
$ws.Range("E18").Value = 5987
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 1607
$ws.Range("B37").Value = 10349
$ws.Range("C37").Value = 536
$ws.Range("D37").Value = 1978
$ws.Range("E37").Value = 8295
$ws.Range("G37").Value = 5
$ws.Range("H37").Value = 76
$ws.Range("B38").Value = 10287
$ws.Range("C38").Value = 929
$ws.Range("D38").Value = 1012
$ws.Range("E38").Value = 9265
$ws.Range("B60").Value = 3075
$ws.Range("C60").Value = 183
$ws.Range("D60").Value = 806
$ws.Range("E60").Value = 2249
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 20
$ws.Range("B61").Value = 2922
$ws.Range("C61").Value = 15
$ws.Range("D61").Value = 2594
$ws.Range("E61").Value = 277
$ws.Range("F61").Value = 61
$ws.Range("H61").Value = 51
$ws.Range("B72").Value = 1645
$ws.Range("C72").Value = 28
$ws.Range("D72").Value = 1139
$ws.Range("E72").Value = 485
$ws.Range("F72").Value = 15
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 21
$ws.Range("B73").Value = 1643
$ws.Range("C73").Value = 8
$ws.Range("D73").Value = 233
$ws.Range("E73").Value = 1361
$ws.Range("F73").Value = 6
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 49
$ws.Range("B80").Value = 1386
$ws.Range("C80").Value = 19
$ws.Range("D80").Value = 500
$ws.Range("E80").Value = 825
$ws.Range("F80").Value = 13
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 61
$ws.Range("B81").Value = 1379
$ws.Range("C81").Value = 6
$ws.Range("D81").Value = 394
$ws.Range("E81").Value = 967
$ws.Range("F81").Value = 5
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 18
$ws.Range("B114").Value = 401
$ws.Range("C114").Value = 21
$ws.Range("D114").Value = 144
$ws.Range("E114").Value = 253
$ws.Range("F114").Value = 4
$ws.Range("H114").Value = 4
$ws.Range("B115").Value = 390
$ws.Range("D115").Value = 8
$ws.Range("E115").Value = 364
$ws.Range("F115").Value = 2
$ws.Range("H115").Value = 18
$ws.Range("D122").Value = 247
$ws.Range("E122").Value = 43
